$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelas")

# --- Row 3 (header row of the "Produto em Stock" table): a new field
#     "Casa" is inserted at B3, pushing the existing headers B3:F3 one
#     column to the right (C3:G3). ---
$ws.Range("F3").Copy($ws.Range("G3"))
$ws.Range("E3").Copy($ws.Range("F3"))
$ws.Range("D3").Copy($ws.Range("E3"))
$ws.Range("C3").Copy($ws.Range("D3"))
$ws.Range("B3").Copy($ws.Range("C3"))
$ws.Range("B3").Value = "Casa"

# --- Row 4 (first data row of that table): a new value (1) for the
#     "Casa" column is inserted at D4, pushing D4:F4 one column to the
#     right (E4:G4). ---
$ws.Range("F4").Copy($ws.Range("G4"))
$ws.Range("E4").Copy($ws.Range("F4"))
$ws.Range("D4").Copy($ws.Range("E4"))
$ws.Range("D4").Value = 1

$ws.Range("F4").Select()
